# Reshape the wide "one row per year" table into a long/melted table:
# a single header row (Unnamed: 0..3) followed by one row per
# category / sub-category with the 2019 / 2018 / 2017 values in B:D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 already carries the bold/centered/bordered header style (style index 1
# in the original file) - reuse it for the rest of the new header row
# B1:D1 before anything else gets clobbered.
$ws.Range("A1").Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Drop the old per-row-is-a-year formatting (A2:A4 used the same bold style)
# and wipe all cell content - the sheet is being fully restructured.
$ws.Cells.ClearContents()
$ws.Range("A2:K4").ClearFormats()

$data = @(
    @("Unnamed: 0", "Unnamed: 1", "Unnamed: 2", "Unnamed: 3"),
    @("(%)", 2019, 2018, 2017),
    @("Flexible working hours", 100, 100, 100),
    @("Full-time", 78.90000000000001, 79.90000000000001, 80.8),
    @("Of which: female", 58.4, 59.9, $null),
    @("Of which: male", 93.90000000000001, 94.59999999999999, $null),
    @("Part-time", 21.1, 20.1, 19.2),
    @("Of which: female", 41.6, 40.1, $null),
    @("Of which: male", 6.1, 5.4, $null),
    @("Virtual offices", 13.2, 12.8, 12.9),
    @("Sabbatical", 0, 0, 0),
    @("Semi-retirement (Altersteilzeit)", 1.76, 1.68, 1.8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $col = $j + 1
        $val = $row[$j]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $col).Value = $val
        }
    }
}

$ws.Range("A1").Select()
